$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.330.50'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '1.877.88'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("D4").Value = "'" + '1.002'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = "'" + '0.7110'
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("D6").Value = "'" + '242.53'
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = "'" + '0.08037'
$ws.Range("E8").Value = '  +3.46%  '
$ws.Range("D9").Value = "'" + '0.3169'
$ws.Range("E9").Value = '  +1.83%  '
$ws.Range("D10").Value = "'" + '25.02'
$ws.Range("E10").Value = '  -0.51%  '
$ws.Range("D11").Value = "'" + '0.08305'
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("D12").Value = '1.886.65'
$ws.Range("E12").Value = '  +0.70%  '
$ws.Range("E13").Value = '  +0.08%  '
$ws.Range("D14").Value = "'" + '94.58'
$ws.Range("E14").Value = '  +3.67%  '
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("D16").Value = "'" + '6.400'
$ws.Range("E16").Value = '  +5.05%  '
$ws.Range("D17").Value = "'" + '0.000008611'
$ws.Range("E17").Value = '  +4.39%  '
$ws.Range("D18").Value = '29.350.33'
$ws.Range("D19").Value = "'" + '242.95'
$ws.Range("E19").Value = '  +0.91%  '
$ws.Range("D20").Value = "'" + '13.31'
$ws.Range("E20").Value = '  +0.47%  '
$ws.Range("D21").Value = '2.138.86'
$ws.Range("E21").Value = '  +0.72%  '
$ws.Range("D22").Value = "'" + '1.002'
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("E23").Value = '  +0.52%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("E25").Value = '  -1.66%  '
$ws.Range("D26").Value = "'" + '9.083'
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("D27").Value = "'" + '162.97'
$ws.Range("E27").Value = '  -0.21%  '
$ws.Range("D28").Value = "'" + '18.54'
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("E29").Value = '  -0.43%  '
$ws.Range("D30").Value = "'" + '4.426'
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").Value = "'" + '4.335'
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("D32").Value = "'" + '1.192'
$ws.Range("E32").Value = '  -7.30%  '
$ws.Range("D33").Value = "'" + '0.05405'
$ws.Range("E33").Value = '  +1.85%  '
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("D35").Value = "'" + '0.7716'
$ws.Range("E35").Value = '  +3.95%  '
$ws.Range("E36").Value = '  +0.67%  '
$ws.Range("D37").Value = "'" + '2.684'
$ws.Range("E37").Value = '  -0.58%  '
$ws.Range("D38").Value = "'" + '0.01890'
$ws.Range("E38").Value = '  +1.04%  '
$ws.Range("D39").Value = '1.266.17'
$ws.Range("E39").Value = '  +3.16%  '
$ws.Range("D40").Value = "'" + '2.753'
$ws.Range("E40").Value = '  +0.80%  '
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("D42").Value = "'" + '113.17'
$ws.Range("E42").Value = '  +2.24%  '
$ws.Range("D43").Value = "'" + '0.9082'
$ws.Range("E43").Value = '  +2.14%  '
$ws.Range("D44").Value = "'" + '74.38'
$ws.Range("E44").Value = '  +1.97%  '
$ws.Range("E45").Value = '  +7.81%  '
$ws.Range("D46").Value = "'" + '1.002'
$ws.Range("E46").Value = '  +0.11%  '
$ws.Range("D47").Value = '2.031.60'
$ws.Range("E47").Value = '  +0.49%  '
$ws.Range("D48").Value = "'" + '1.807'
$ws.Range("E48").Value = '  -0.39%  '
$ws.Range("D49").Value = "'" + '0.5225'
$ws.Range("E49").Value = '  +0.21%  '
$ws.Range("D50").Value = "'" + '9.473'
$ws.Range("E50").Value = '  +0.26%  '
$ws.Range("D51").Value = "'" + '0.4375'
$ws.Range("E51").Value = '  +1.23%  '
